# T06 Scrum Planning - "Added Prioroty form PO"
# Adds a "Priority" column (D) with Hight/Mid/Low values, a notes column (I)
# with prerequisite text, a couple of "1. Sprint" markers in column L, and a
# new backlog item (row 9) on the "Product Backlog" sheet. Also switches the
# active/selected sheet from "ProjectTeam" to "Product Backlog".

$wb = $excel.ActiveWorkbook

$ws    = $wb.Worksheets.Item("Product Backlog")
$refWs = $wb.Worksheets.Item("Sprint Backlog")

# ---------------------------------------------------------------------
# 1. New row 8 (id 8) - "Optisches Grundgerüst / Haupt UI"
# ---------------------------------------------------------------------
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Optisches Grundgerüst / Haupt UI"

# ---------------------------------------------------------------------
# 2. Priority column (D) values for every backlog item
# ---------------------------------------------------------------------
$ws.Range("D2").Value = "Hight"
$ws.Range("D3").Value = "Low"
$ws.Range("D4").Value = "Mid"
$ws.Range("D5").Value = "Hight"
$ws.Range("D6").Value = "Hight"
$ws.Range("D7").Value = "Low"
$ws.Range("D8").Value = "Low"
$ws.Range("D9").Value = "Hight"

# ---------------------------------------------------------------------
# 3. New "prerequisite" notes column (I) + a couple of sprint markers (L)
# ---------------------------------------------------------------------
$ws.Range("I2").Value = "Voraussetzung Patientenakte"
$ws.Range("I5").Value = "Voraussetzung für Notzien"
$ws.Range("I6").Value = "Im ersten Sprint evt. Spliten"
$ws.Range("I7").Value = "Superimplemetation: Medikamente"
$ws.Range("I8").Value = "Superimplemetation: Medikamente"

$ws.Range("L6").Value = "1. Sprint"
$ws.Range("L9").Value = "1. Sprint"

# ---------------------------------------------------------------------
# 4. Formatting - reuse existing cell formats instead of inventing new
#    fonts/styles, by copying format only (keeps style table minimal,
#    same as what Excel itself would normally collapse to).
# ---------------------------------------------------------------------

# Column A (ID) vertical-top alignment, same as A2, for the new/older rows
$ws.Range("A2").Copy()
$ws.Range("A3:A9").PasteSpecial(-4122)

# Column D (Priority) uses the same vertical-top alignment as column A
$ws.Range("A2").Copy()
$ws.Range("D2:D9").PasteSpecial(-4122)

# Column B (Story Name) - new row 9 should look like the rest of column B
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)

# Column C (Description): switch from "wrap only" to "top + wrap" (same
# format already used for similar description cells on the Sprint Backlog
# sheet), for every existing description row 2-7
$refWs.Range("D2").Copy()
$ws.Range("C2:C7").PasteSpecial(-4122)

# Column I (new notes column): same "top + wrap" alignment
$refWs.Range("D2").Copy()
$ws.Range("I2:I8").PasteSpecial(-4122)

# C8 keeps its own (bold-ish) font/left-align/wrap format, just switch the
# vertical alignment from centered to top, matching the rest of the sheet.
# (-4160 == xlTop)
$ws.Range("C8").VerticalAlignment = -4160

$excel.CutCopyMode = 0

# Column I needs to be wide enough for the new notes text
$ws.Columns.Item(9).ColumnWidth = 30.140625

# ---------------------------------------------------------------------
# 5. Selection / active sheet bookkeeping - "Product Backlog" becomes the
#    active tab (previously "ProjectTeam" was selected).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("C6").Select()
